$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.422
$ws.Range("D5").Value = 0.526
$ws.Range("E5").Value = 0.551
$ws.Range("F5").Value = 0.582
$ws.Range("G5").Value = 0.514
$ws.Range("H5").Value = 0.536

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.422
$ws.Range("D7").Value = 0.526

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.432
$ws.Range("D8").Value = 0.602
$ws.Range("E8").Value = 0.627
$ws.Range("F8").Value = 0.666
$ws.Range("G8").Value = 0.625
$ws.Range("H8").Value = 0.651

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.437
$ws.Range("C9").Value = 0.5580000000000001
$ws.Range("D9").Value = 0.658
$ws.Range("E9").Value = 0.676
$ws.Range("F9").Value = 0.707
$ws.Range("G9").Value = 0.642
$ws.Range("H9").Value = 0.66
